# corrected data cleaning for pre/post/total fixation data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Strip the header styling (bold font + border) from row 1 and clear the
#    "Unnamed: 0" label out of A1, leaving it blank like the rest of the row.
# ---------------------------------------------------------------------------
$ws.Range("A1:M1").ClearFormats()
$ws.Range("A1").Value = ""

# ---------------------------------------------------------------------------
# 2. Update the corrected values for the first (pre) fixation-metrics block,
#    rows 3-7 (row 8 is untouched by the correction).
# ---------------------------------------------------------------------------
$updates = @{
    "B3" = 21;        "C3" = 0;       "D3" = 9;        "G3" = 18;       "I3" = 24
    "B4" = 55;        "C4" = 1;       "D4" = 25;       "G4" = 44;       "I4" = 201
    "B5" = 17932.72;  "C5" = 183.48;  "D5" = 10160.3;  "G5" = 14931.83; "I5" = 81543.5
    "B6" = 11.32;     "C6" = 0.12;    "D6" = 6.42;     "F6" = 0.11;     "G6" = 9.43;  "I6" = 51.49; "K6" = 0.23
    "B7" = 326.05;    "C7" = 183.48;  "D7" = 406.41;   "G7" = 339.36;   "I7" = 405.69
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# ---------------------------------------------------------------------------
# 3. Turn the previously-blank row 11 into a repeated header row (same labels
#    as row 1, but column A stays blank and none of these cells are bold /
#    bordered since that styling was removed workbook-wide).
# ---------------------------------------------------------------------------
$headers = @{
    "B" = "code"; "C" = "declaration"; "D" = "exception"; "E" = "gemini"
    "F" = "index"; "G" = "loop body";  "H" = "param";     "I" = "sum"
    "J" = "var";   "K" = "var2";       "L" = "var3";      "M" = "var4"
}
foreach ($col in $headers.Keys) {
    $ws.Range($col + "11").Value = $headers[$col]
}

# ---------------------------------------------------------------------------
# 4. Append a brand-new "post" fixation-metrics block in rows 12-18, mirroring
#    the structure of rows 2-8.
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = "Fixation based metrics"

$labels = @{
    13 = "Revisit count"
    14 = "Fixation count"
    15 = "Dwell time (ms)"
    16 = "Dwell time (%)"
    17 = "Fixation duration (ms)"
    18 = "First fixation duration (ms)"
}
foreach ($r in $labels.Keys) {
    $ws.Range("A" + $r).Value = $labels[$r]
}

$newData = @{
    "B13" = 8;        "C13" = 0;       "D13" = 0;       "E13" = 23;       "G13" = 7;       "I13" = 27
    "B14" = 10;       "C14" = 1;       "D14" = 1;       "E14" = 177;      "G14" = 8;       "I14" = 133
    "B15" = 3503.83;  "C15" = 200.17;  "D15" = 417.27;  "E15" = 41725.81; "G15" = 3220.22; "I15" = 41804.07
    "B16" = 2.35;     "C16" = 0.13;    "D16" = 0.28;    "E16" = 28;       "G16" = 2.16;    "I16" = 28.06
    "B17" = 350.38;   "C17" = 200.17;  "D17" = 417.27;  "E17" = 235.74;   "G17" = 402.53;  "I17" = 314.32
    "B18" = 417.27;   "C18" = 200.17;  "D18" = 417.27;  "E18" = 333.44;   "G18" = 417.27;  "I18" = 199.94
}
foreach ($addr in $newData.Keys) {
    $ws.Range($addr).Value = $newData[$addr]
}
